$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.8330250000000001
$ws.Range("H2").Value = 2.499075
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.668351
$ws.Range("N2").Value = 2.005053
$ws.Range("O2").Value = 0.04158981742241631
$ws.Range("P2").Value = 0.0415898174224163
$ws.Range("Q2").Value = 0.5567530917750001
$ws.Range("R2").Value = 5.010777825975001
$ws.Range("S2").Value = 0.04158981742241631
$ws.Range("T2").Value = 0.0415898174224163

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.8330250000000001
$ws.Range("H3").Value = 2.499075
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.435983
$ws.Range("N3").Value = 31.307949
$ws.Range("O3").Value = 0.649405219104094
$ws.Range("P3").Value = 0.649405219104094
$ws.Range("Q3").Value = 8.693434738575
$ws.Range("R3").Value = 78.24091264717501
$ws.Range("S3").Value = 0.649405219104094
$ws.Range("T3").Value = 0.649405219104094

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.8330250000000001
$ws.Range("H4").Value = 2.499075
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.965729333333333
$ws.Range("N4").Value = 14.897188
$ws.Range("O4").Value = 0.3090049634734898
$ws.Range("P4").Value = 0.3090049634734898
$ws.Range("Q4").Value = 4.1365766779
$ws.Range("R4").Value = 37.2291901011
$ws.Range("S4").Value = 0.3090049634734898
$ws.Range("T4").Value = 0.3090049634734898
